$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 378.8
$ws.Range("J103").Value = 500
$ws.Range("L103").Value = 1500
$ws.Range("N103").Value = -2672
$ws.Range("H113").Value = 12827.158
$ws.Range("I113").Value = 12645.083
$ws.Range("J113").Value = 13139.286
$ws.Range("K113").Value = 12645.083
$ws.Range("L113").Value = 13139.286
$ws.Range("M113").Value = -9391.083000000001
$ws.Range("N113").Value = -19647.286
$ws.Range("H125").Value = 3172.8823
$ws.Range("J125").Value = 2964
$ws.Range("L125").Value = 26676
$ws.Range("N125").Value = -31596
$ws.Range("H132").Value = 1205.0139
$ws.Range("I132").Value = 1193.194
$ws.Range("K132").Value = 3579.582
$ws.Range("M132").Value = -1049.582
$ws.Range("H137").Value = 2899.8386
$ws.Range("I137").Value = 3028.7368
$ws.Range("J137").Value = 2695.75
$ws.Range("K137").Value = 9086.2104
$ws.Range("L137").Value = 8087.25
$ws.Range("M137").Value = -6536.2104
$ws.Range("N137").Value = -13187.25
$ws.Range("H138").Value = 4313.95
$ws.Range("I138").Value = 2059.0588
$ws.Range("J138").Value = 4775.7954
$ws.Range("K138").Value = 6177.176399999999
$ws.Range("L138").Value = 14327.3862
$ws.Range("M138").Value = -1037.176399999999
$ws.Range("N138").Value = -24607.3862
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10982.197
$ws.Range("I32").Value = 4585.7363
$ws.Range("J32").Value = 43878.285
$ws.Range("K32").Value = 4585.7363
$ws.Range("L32").Value = 43878.285
$ws.Range("M32").Value = -4298.7363
$ws.Range("N32").Value = -44452.285
$ws.Range("H45").Value = 24713.572
$ws.Range("I45").Value = 41374
$ws.Range("K45").Value = 41374
$ws.Range("M45").Value = -40997
$ws.Range("H61").Value = 198229.42
$ws.Range("I61").Value = 2286.6365
$ws.Range("K61").Value = 2286.6365
$ws.Range("M61").Value = -2074.6365
$ws.Range("H97").Value = 1810.6
$ws.Range("I97").Value = 2211.1428
$ws.Range("J97").Value = 876
$ws.Range("K97").Value = 2211.1428
$ws.Range("L97").Value = 876
$ws.Range("M97").Value = -1715.1428
$ws.Range("N97").Value = -1868
$ws.Range("H102").Value = 5873
$ws.Range("J102").Value = 5500
$ws.Range("L102").Value = 5500
$ws.Range("N102").Value = -8744
$ws.Range("H136").Value = 198229.42
$ws.Range("I136").Value = 2286.6365
$ws.Range("K136").Value = 6859.9095
$ws.Range("M136").Value = -4309.9095
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 29998
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 29998
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 29998
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -30992
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H99").Value = 4733.0415
$ws.Range("I99").Value = 3176.3845
$ws.Range("J99").Value = 6572.727
$ws.Range("K99").Value = 3176.3845
$ws.Range("L99").Value = 6572.727
$ws.Range("M99").Value = -1678.3845
$ws.Range("N99").Value = -9568.726999999999
$ws.Range("H105").Value = 3203.6428
$ws.Range("I105").Value = 3203.6428
$ws.Range("K105").Value = 3203.6428
$ws.Range("M105").Value = -1456.6428
$ws.Range("H107").Value = 1752.1111
$ws.Range("I107").Value = 1791.8077
$ws.Range("J107").Value = 720
$ws.Range("K107").Value = 1791.8077
$ws.Range("L107").Value = 720
$ws.Range("M107").Value = 128.1922999999999
$ws.Range("N107").Value = -4560
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 1699.8
$ws.Range("I25").Value = 1699.8
$ws.Range("K25").Value = 1699.8
$ws.Range("M25").Value = -1525.8
$ws.Range("H31").Value = 49918
$ws.Range("I31").Value = 60414.59
$ws.Range("K31").Value = 60414.59
$ws.Range("M31").Value = -60119.59
$ws.Range("H34").Value = 49918
$ws.Range("I34").Value = 60414.59
$ws.Range("K34").Value = 60414.59
$ws.Range("M34").Value = -60212.59
$ws.Range("H88").Value = 13795.111
$ws.Range("I88").Value = 17286.5
$ws.Range("J88").Value = 12797.571
$ws.Range("K88").Value = 17286.5
$ws.Range("L88").Value = 12797.571
$ws.Range("M88").Value = -16880.5
$ws.Range("N88").Value = -13609.571
$ws.Range("H91").Value = 13795.111
$ws.Range("I91").Value = 17286.5
$ws.Range("J91").Value = 12797.571
$ws.Range("K91").Value = 17286.5
$ws.Range("L91").Value = 12797.571
$ws.Range("M91").Value = -15882.5
$ws.Range("N91").Value = -15605.571
$ws.Range("H107").Value = 852.3200000000001
$ws.Range("I107").Value = 562.3684
$ws.Range("K107").Value = 562.3684
$ws.Range("M107").Value = 1357.6316
$ws.Range("H134").Value = 11371.5
$ws.Range("I134").Value = 6166.3184
$ws.Range("K134").Value = 18498.9552
$ws.Range("M134").Value = -15963.9552
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 83352.664
$ws.Range("I9").Value = 29
$ws.Range("K9").Value = 87
$ws.Range("M9").Value = 137
$ws.Range("H44").Value = 2810.4
$ws.Range("I44").Value = 1660
$ws.Range("J44").Value = 3960.8
$ws.Range("K44").Value = 4980
$ws.Range("L44").Value = 11882.4
$ws.Range("M44").Value = -4582
$ws.Range("N44").Value = -12678.4
$ws.Range("H54").Value = 1000
$ws.Range("J54").Value = 1000
$ws.Range("L54").Value = 3000
$ws.Range("N54").Value = -4118
$ws.Range("H81").Value = 3726.6667
$ws.Range("H84").Value = 3726.6667
$ws.Range("H128").Value = 203497
$ws.Range("I128").Value = 203497
$ws.Range("K128").Value = 610491
$ws.Range("M128").Value = -605511
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2165.4119
$ws.Range("I107").Value = 1104.2727
$ws.Range("K107").Value = 1104.2727
$ws.Range("M107").Value = 815.7273
$ws.Range("H122").Value = 2966.25
$ws.Range("I122").Value = 2751.3333
$ws.Range("K122").Value = 8253.999899999999
$ws.Range("M122").Value = -5803.999899999999
$ws.Range("H126").Value = 19046.428
$ws.Range("I126").Value = 23104.545
$ws.Range("J126").Value = 4166.6665
$ws.Range("K126").Value = 69313.63499999999
$ws.Range("L126").Value = 12499.9995
$ws.Range("M126").Value = -66843.63499999999
$ws.Range("N126").Value = -17439.9995
$ws.Range("H133").Value = 103046.5
$ws.Range("J133").Value = 103046.5
$ws.Range("L133").Value = 103046.5
$ws.Range("N133").Value = -113166.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4182.9546
$ws.Range("I7").Value = 3762.7693
$ws.Range("K7").Value = 3762.7693
$ws.Range("M7").Value = -3650.7693
$ws.Range("H22").Value = 2772.5173
$ws.Range("I22").Value = 1908.2307
$ws.Range("J22").Value = 3474.75
$ws.Range("K22").Value = 1908.2307
$ws.Range("L22").Value = 3474.75
$ws.Range("M22").Value = -1613.2307
$ws.Range("N22").Value = -4064.75
$ws.Range("H27").Value = 2772.5173
$ws.Range("I27").Value = 1908.2307
$ws.Range("J27").Value = 3474.75
$ws.Range("K27").Value = 1908.2307
$ws.Range("L27").Value = 3474.75
$ws.Range("M27").Value = -1801.2307
$ws.Range("N27").Value = -3688.75
$ws.Range("H40").Value = 4202.237
$ws.Range("I40").Value = 3541
$ws.Range("J40").Value = 6332.8887
$ws.Range("K40").Value = 3541
$ws.Range("L40").Value = 6332.8887
$ws.Range("M40").Value = -3405
$ws.Range("N40").Value = -6604.8887
$ws.Range("H46").Value = 3899
$ws.Range("J46").Value = 4198.6665
$ws.Range("L46").Value = 4198.6665
$ws.Range("N46").Value = -4574.6665
$ws.Range("H61").Value = 2141.3958
$ws.Range("I61").Value = 1646.2094
$ws.Range("J61").Value = 6400
$ws.Range("K61").Value = 1646.2094
$ws.Range("L61").Value = 6400
$ws.Range("M61").Value = -1444.2094
$ws.Range("N61").Value = -6804
$ws.Range("H68").Value = 2974.96
$ws.Range("I68").Value = 3167.4119
$ws.Range("K68").Value = 3167.4119
$ws.Range("M68").Value = -2418.4119
$ws.Range("H71").Value = 2974.96
$ws.Range("I71").Value = 3167.4119
$ws.Range("K71").Value = 15837.0595
$ws.Range("M71").Value = -12093.0595
$ws.Range("H113").Value = 2141.3958
$ws.Range("I113").Value = 1646.2094
$ws.Range("J113").Value = 6400
$ws.Range("K113").Value = 1646.2094
$ws.Range("L113").Value = 6400
$ws.Range("M113").Value = 523.7906
$ws.Range("N113").Value = -10740
$ws.Range("H126").Value = 4182.9546
$ws.Range("I126").Value = 3762.7693
$ws.Range("K126").Value = 11288.3079
$ws.Range("M126").Value = -8818.3079
$ws.Range("H132").Value = 5424.5454
$ws.Range("I132").Value = 6000
$ws.Range("J132").Value = 5208.75
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 15626.25
$ws.Range("M132").Value = -15470
$ws.Range("N132").Value = -20686.25
$ws.Range("H136").Value = 7262.375
$ws.Range("I136").Value = 6800
$ws.Range("K136").Value = 20400
$ws.Range("M136").Value = -17850
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12797.2
$ws.Range("I62").Value = 6332
$ws.Range("J62").Value = 22495
$ws.Range("K62").Value = 6332
$ws.Range("L62").Value = 22495
$ws.Range("M62").Value = -5708
$ws.Range("N62").Value = -23743
$ws.Range("H65").Value = 12797.2
$ws.Range("I65").Value = 6332
$ws.Range("J65").Value = 22495
$ws.Range("K65").Value = 31660
$ws.Range("L65").Value = 112475
$ws.Range("M65").Value = -28540
$ws.Range("N65").Value = -118715
$ws.Range("H107").Value = 1173.7142
$ws.Range("I107").Value = 1178.0834
$ws.Range("K107").Value = 3534.2502
$ws.Range("M107").Value = -1614.2502
$ws.Range("H113").Value = 334818.56
$ws.Range("I113").Value = 1241.826
$ws.Range("K113").Value = 3725.478
$ws.Range("M113").Value = -1555.478
